# Apply the "index" field addition to the generated JS object literals in
# column D (rows 2-26), mirroring the author's change to the formulas:
#   ... htmlContent: '' },            ->   ... htmlContent: '', index: <A-1> },
#
# Each row's formula is re-entered through the COM object model so Excel
# recalculates the shared-formula group (D3:D26) and refreshes every
# cached <v> value automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 26; $row++) {
    $formula = ' "{ x:" & B' + $row + ' & ", z:" & C' + $row + ' & ", id: ''p" & A' + $row + ' & "'', imageUrl: ''assets/" & A' + $row + ' & "'', vaov: 180, vOffset: 0, maxpitch: 40, minpitch: -90, htmlContent: '''', index: " & A' + $row + ' -1 &" },"'
    $ws.Range("D$row").Formula = "=$formula"
}
